$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.216.51'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +3.39%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.605.79'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.03%  '
$ws.Range('E4').Value = '  -0.49%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.79'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.23%  '
$ws.Range('E6').Value = '  -0.50%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.485'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.41%  '
$ws.Range('E8').Value = '  +1.95%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0618'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.47%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.06'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.53%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0816'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.19%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.829.61'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.16%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.602.24'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.70%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.02'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.60%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.511'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.83%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '26.187.40'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '60.58'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.15%  '
$ws.Range('E18').Value = '  +2.10%  '
$ws.Range('E19').Value = '  -0.65%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '197.80'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +6.31%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.25'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.39%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.40'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.54%  '
$ws.Range('E23').Value = '  +1.72%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '142.40'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.77%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.76'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.48%  '
$ws.Range('B26').Value = 'BinanceUSD'
$ws.Range('C26').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.61%  '
$ws.Range('B27').Value = 'Stellar'
$ws.Range('C27').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.126'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.32%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.15'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.60%  '
$ws.Range('E29').Value = '  -0.21%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.18'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.19%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0471'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.05%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.14'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.39%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.01'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.16%  '
$ws.Range('E34').Value = '  +0.87%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.36'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.43%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.107.22'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.84%  '
$ws.Range('E37').Value = '  -0.17%  '
$ws.Range('E38').Value = '  +1.03%  '
$ws.Range('E39').Value = '  +0.31%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.786'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.57%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.499'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.95%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.775'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.65%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.741.54'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.11%  '
$ws.Range('E44').Value = '  +1.25%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '92.43'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.14%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0₆0108'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.14%  '
$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.55'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +9.51%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '53.45'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.98%  '
$ws.Range('E49').Value = '  -0.16%  '
$ws.Range('E50').Value = '  +0.41%  '
$ws.Range('E51').Value = '  -0.15%  '
